$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line rows ("line7", "line8") are inserted right after the
# existing "line6" row, which pushes the old contingency rows
# (extr1..extr8, previously rows 8-15) down to rows 10-17. Their
# from_bus / to_bus / in_service values are also recomputed.
#
# columns: A = index, B = name, C = from_bus, D = to_bus, E = in_service
$data = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

# Rows 16 and 17 are brand new (the sheet previously only went down to
# row 15), so column A needs the same bold/centered/bordered style that
# is already used by the other index cells in column A.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
}
